# Update the "timestamp" column (O) for all data rows (2-64) to the new
# crawl timestamp, reflecting a re-upload of the scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-09-06 20:59:47"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 64
}

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
